{"js": "// Replace the 25 division-fact cells in the single table, in document\n// order, with their new values (per the commit diff). Several source\n// strings repeat (e.g. \"55\u00f78=\", \"40\u00f75=\", \"98\u00f79=\") so we must walk the\n// table row-by-row / cell-by-cell rather than do a single global\n// find/replace of the source text.\nconst newValues = [\n  [\"77\u00f75=\", \"35\u00f78=\", \"59\u00f72=\", \"40\u00f79=\", \"96\u00f78=\"],\n  [\"77\u00f74=\", \"40\u00f72=\", \"80\u00f73=\", \"71\u00f72=\", \"25\u00f75=\"],\n  [\"47\u00f79=\", \"63\u00f72=\", \"58\u00f79=\", \"41\u00f79=\", \"63\u00f76=\"],\n  [\"31\u00f73=\", \"52\u00f79=\", \"37\u00f74=\", \"61\u00f79=\", \"88\u00f78=\"],\n  [\"30\u00f76=\", \"16\u00f79=\", \"33\u00f74=\", \"90\u00f74=\", \"46\u00f74=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Only every 4th row (0, 4, 8, 12, 16) holds equation text; the rows in\n// between are blank \"answer\" rows and are left untouched.\nfor (let r = 0; r < newValues.length; r++) {\n  const rowIndex = r * 4;\n  const cells = rows.items[rowIndex].cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < newValues[r].length; c++) {\n    const paras = cells.items[c].body.paragraphs;\n    paras.load(\"items\");\n    await context.sync();\n\n    // insertText(\"Replace\") on the Paragraph keeps the existing\n    // paragraph/run formatting (font, size, alignment) intact, unlike\n    // calling it on the TableCell.body, which drops rPr/pPr.\n    paras.items[0].insertText(newValues[r][c], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-fact cells in the single table, in document\n# order, with their new values (per the commit diff). Several source\n# strings repeat (e.g. \"55\u00f78=\", \"40\u00f75=\", \"98\u00f79=\") so we address cells\n# positionally (Table.Cell(row, col)) rather than doing a single global\n# Find/Replace of the source text.\n$newValues = @(\n    @(\"77\u00f75=\", \"35\u00f78=\", \"59\u00f72=\", \"40\u00f79=\", \"96\u00f78=\"),\n    @(\"77\u00f74=\", \"40\u00f72=\", \"80\u00f73=\", \"71\u00f72=\", \"25\u00f75=\"),\n    @(\"47\u00f79=\", \"63\u00f72=\", \"58\u00f79=\", \"41\u00f79=\", \"63\u00f76=\"),\n    @(\"31\u00f73=\", \"52\u00f79=\", \"37\u00f74=\", \"61\u00f79=\", \"88\u00f78=\"),\n    @(\"30\u00f76=\", \"16\u00f79=\", \"33\u00f74=\", \"90\u00f74=\", \"46\u00f74=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Only every 4th row (1, 5, 9, 13, 17 in Word's 1-based row numbering)\n# holds equation text; the rows in between are blank \"answer\" rows and\n# are left untouched.\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $rowIndex = 1 + $r * 4\n    $rowValues = $newValues[$r]\n    for ($c = 0; $c -lt $rowValues.Length; $c++) {\n        $cell = $t.Cell($rowIndex, $c + 1)\n        # Assigning Range.Text replaces the run's text in place, keeping\n        # the existing run/paragraph formatting (font, size, alignment).\n        $cell.Range.Text = $rowValues[$c]\n    }\n}\n"}
